# Evidence workbook update: mint NFTs on IRISnet (sheet "A2") and
# record the per-NFT transaction hashes + minted NFT ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A2")

# Make "A2" the active sheet/tab (workbook activeTab moves from "A1" to "A2").
$ws.Activate()

$txHashes = @(
    "E424A5239F1FB5B5DF698111AD8A993F9391B8CC23997281FE65182137218C5D",
    "03FB72A529F5CB9C3FCF7511BBAD9E87E46DA14BEA046B39FF79CD32427BF032",
    "9E5879D98778E8777253AF21809B38CD5A17E1ACD3BD618DBACE389294ECAA8A",
    "4B2ADC13B0FFF9FF5FE2C795D89FE6DF3E49D32423B87F2F4EA216ACE1C0B865",
    "AB38CDFA336E84839D66F021006268F74DEF6FE04388DEFEFB16968282D0039A",
    "AC2C856CD9467A77111E999A249BCB07BD377A59A6539CC186E5E38940C2FAC8"
)

$classId = "arkprotocol002"

$nftIds = @(
    "arkNFT001",
    "arkNFT002",
    "arkNFT003",
    "arkNFT004",
    "arkNFT005",
    "arkNFT006"
)

for ($i = 0; $i -lt $txHashes.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $txHashes[$i]
    $ws.Cells.Item($row, 2).Value = $classId
    $ws.Cells.Item($row, 3).Value = $nftIds[$i]
}

# A few trailing formatted-but-empty cells below the data (carried over from
# the original sheet's used range / formatting extent).
$ws.Range("B8").Style = $ws.Range("B7").Style
$ws.Range("B9").Style = $ws.Range("B7").Style
$ws.Range("B10").Style = $ws.Range("B7").Style

$ws.Range("B13").Select()
